$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 100,4
$arr[0,0] = 'I recently upgraded to the Xiaomi X Series 4K UHD Google TV, and I’m quite impressed by what it delivers at its price point. Here’s my personal take:

Picture Quality:

The 4K resolution with Dolby Vision really stands out. Watching UHD content is a treat, with vibrant colors and clear details, especially in brighter scenes. However, when streaming SDR content or older shows, you might notice some graininess, especially on larger screens like the 65-inch model. Adjusting the backlight and saturation can help if you prefer more natural tones, as the default settings can be a bit too vivid for my taste.

Sound:

The 30W speakers are surprisingly good for regular TV viewing, supporting Dolby Audio and DTS. For everyday use, they’re loud and clear enough, but I found the bass to be lacking, so if you''re a fan of deep, immersive sound, investing in a soundbar might be a good idea.

User Interface:

Google TV is easy to navigate, and I love how it suggests shows and movies across platforms based on my viewing habits. Plus, Xiaomi’s PatchWall feature is a great addition, giving access to over 200 live channels without extra subscriptions. Switching between apps is generally smooth, though I noticed slight lag after prolonged use, but it’s nothing that impacts the overall experience too much.

Design:

The sleek black design with thin bezels looks great in any room. I opted to mount mine on the wall, which I think suits the larger models better as they can take up a lot of space on a stand. It''s not a premium build, but it’s sturdy and fits well with modern interiors.

Smoothness
In terms of smoothness, the Xiaomi X Series 4K UHD Google TV is decent but not flawless. The TV is powered by a quad-core A55 chipset, paired with 2GB of RAM and 8GB of internal storage. For regular tasks like switching between apps, browsing, and streaming, the performance is generally smooth. However, some users have reported occasional lag, especially after prolonged use or when loading heavier apps. It''s not a dealbreaker, but the experience might not be as snappy as higher-end models.

The Google TV interface is a big plus here, as it runs efficiently, with intuitive navigation and quick access to content across various platforms. Xiaomi’s PatchWall feature adds to the smoothness, with seamless integration of live TV and streaming apps into one unified experience.

For basic daily use, the smoothness is acceptable, but if you’re after a flawless, zero-lag experience, you might notice occasional slowdowns with heavy multitasking.

Final Thoughts:

Overall, the Xiaomi X Series delivers excellent value for anyone seeking a 4K TV with smart features at an affordable price. Its picture quality shines with UHD content, and the Google TV integration is intuitive. Just be prepared for some occasional performance hiccups and consider adding external speakers if audio quality is a priority for you.'
$arr[0,1] = 'Style Name: Google X Series 2023Size: 43'
$arr[0,2] = '5.0 out of 5 stars'
$arr[0,3] = 'Placeholder'
$arr[1,0] = 'Very good and colourful panel, deep blacks, moderate viewing angles. Metal back feels solid. Speaker is also lound and clear.
One con is that TV sometimes get slow and ui frames caps to like 30 hz, but it is rare. Hope this issue gets solved by update.
Overall if anyone wants fabulous picture quality this should be a very good option.
I got the tv at 21750 during sale, so it is value for money product'
$arr[1,1] = 'Style Name: A Series 2024Size: 43 inches'
$arr[1,2] = '4.0 out of 5 stars'
$arr[1,3] = $null
$arr[2,0] = 'I have purchased this MI tv on 26th January. Amazon delivered timely but installation was late due to local player. I was raised the installation complain to MI customer care.They were responded well and ensured the job. After two months usage , I feel happy for it. Sometimes it lags during streaming,Voice command response slow, screencast through mobile not happening, contrast is high , visible in black hair with white light. Picture and sound quality is nice. I m using Jio fiber and observed sound variation on different applications like jio,Zee5,sony. Overall it''s good product. Life of LED is 4 to 7 years . It depends on usage. Recommend this unit for purchase.??'
$arr[2,1] = 'Style Name: A Series 2024Size: 43 inches'
$arr[2,2] = '5.0 out of 5 stars'
$arr[2,3] = 'ΞυηκηοωηΞ'
$arr[3,0] = 'Brightness and Video Quality is Best. 4k UHD Videos are Crisp to view. Dolby audio only available in Optical Format and not in HDMI ARC. Audio output is also Good Quality. Overall Best TV. Installation done by Amazon person on delivery day itself quickly and efficiently which is a plus point.'
$arr[3,1] = 'Style Name: A Series 2024Size: 43 inches'
$arr[3,2] = '4.0 out of 5 stars'
$arr[3,3] = $null
$arr[4,0] = 'Best tv in the segment. I bought this tv after a long research in YouTube and many articles.

Pros :
1. Best tv in the price range , I bought it under 30k by using card offers
2. The performance of the tv is really well and it''s really bright and watching 4k video was awesome
3. The sound of tv is super, no need of sub woofer until you need a theatre experience.
4. The viewing angle was good , as I tried from few side angles

Cons :
1. After a month , I felt some lag in the tv perfomance, but after a watching a youtube video , I reduced the settings and animation and few other things , that really boosted the speed of the tv
------------

I raised a service request and they also responded soon , but we can''t a expect a speed of what we have in our mobile phones, because it''s of minimum ram right?
I haven''t used the storage option to store my photos , as it may slow down the performance.

Overall it''s a best tv in the price segment good to buy'
$arr[4,1] = 'Style Name: Google X Series 2023Size: 50'
$arr[4,2] = '5.0 out of 5 stars'
$arr[4,3] = 'Santoshkumar G. Patel'
$arr[5,0] = 'I recently purchased the MI 108 cm (43 inches) X Series 4K Ultra HD Smart Google TV (L43M8-A2IN), and it has completely transformed my viewing experience! Here’s my in-depth review of this impressive smart TV.

Design and Build Quality
The design of the MI X Series is sleek and modern, with a minimalistic bezel that maximizes the screen space. The black finish gives it a sophisticated look that blends well with any home decor. The build quality feels solid, and the lightweight design makes it easy to mount on a wall or place on a stand.

Picture Quality
The 4K Ultra HD resolution is stunning. The picture clarity, vibrant colors, and sharp details make watching movies and shows a joy. Whether I’m streaming in high definition or watching regular broadcasts, the upscaling technology does an excellent job of enhancing image quality. The viewing angles are impressive too, allowing for a great experience from different parts of the room.

Smart Features
The integration of Google TV is a major highlight. Navigating through the interface is intuitive, and I love the personalized recommendations based on my viewing habits. The voice search feature works seamlessly, making it easy to find content without scrolling through endless menus. Additionally, access to a wide range of apps, including Netflix, YouTube, and Prime Video, ensures I never run out of things to watch.

Sound Quality
The audio quality is surprisingly good for a TV of this size. The sound is clear, and while it may not replace a dedicated sound system, it’s adequate for regular viewing. For those who prefer a more immersive experience, connecting external speakers or a soundbar is a breeze.

Connectivity
The TV comes with multiple HDMI and USB ports, allowing for easy connections to gaming consoles, Blu-ray players, and other devices. The Wi-Fi connectivity is stable, ensuring smooth streaming without interruptions. I also appreciate the Bluetooth capability for connecting wireless headphones and speakers.

Remote Control
The remote is user-friendly, featuring dedicated buttons for popular streaming services. The Google Assistant integration is a fantastic addition, allowing for voice commands to control the TV and compatible smart home devices.

Conclusion
Overall, the MI 108 cm (43 inches) X Series 4K Ultra HD Smart Google TV L43M8-A2IN is an excellent choice for anyone looking to upgrade their home entertainment system. With its stunning picture quality, smart features, and solid build, it offers great value for the price. Whether you’re a movie buff, gamer, or casual viewer, this TV will enhance your viewing experience significantly. I highly recommend it!'
$arr[5,1] = 'Style Name: Google X Series 2023Size: 43'
$arr[5,2] = '5.0 out of 5 stars'
$arr[5,3] = $null
$arr[6,0] = 'There is no dolby audio in this TV. I can see it was mentioned on product description but in reality no dolby audio. Picture quality is good but, poor audio clarity with no bass. NO DOLBY Audio.'
$arr[6,1] = 'Style Name: A Series 2024Size: 43 inches'
$arr[6,2] = '3.0 out of 5 stars'
$arr[6,3] = 'G. B.'
$arr[7,0] = 'Very good sound quality and display is nice
Budget fit'
$arr[7,1] = 'Style Name: Google X Android SeriesSize: 50 inches'
$arr[7,2] = '4.0 out of 5 stars'
$arr[7,3] = 'kamal Bharathi'
$arr[8,0] = 'Good'
$arr[8,1] = 'Style Name: Google X Series 2024Size: 50'
$arr[8,2] = '4.0 out of 5 stars'
$arr[8,3] = 'Shashi Kant'
$arr[9,0] = 'Very good ?? ?? percent'
$arr[9,1] = 'Style Name: A Series 2024Size: 43 inches'
$arr[9,2] = '4.0 out of 5 stars'
$arr[9,3] = 'Amazon Customer'
$arr[10,0] = 'Good'
$arr[10,1] = 'Style Name: A Series 2024Size: 43 inches'
$arr[10,2] = '4.0 out of 5 stars'
$arr[10,3] = 'Sushant padwal'
$arr[11,0] = 'After using 2 years I am writing this review. Worth to buy.
Superb performance and value for money ??
Installation done by Service provider.
Picture quality is still same.
Easy to operate by remote control because it doesn’t have lot of buttons.'
$arr[11,1] = 'Style Name: Google X Android SeriesSize: 55 inches'
$arr[11,2] = '4.0 out of 5 stars'
$arr[11,3] = 'Jeffvin'
$arr[12,0] = 'Click to play video
Best and best TV and sound quality
Delevery boy charge 500 rs'
$arr[12,1] = 'Style Name: Google X Series 2024Size: 50'
$arr[12,2] = '4.0 out of 5 stars'
$arr[12,3] = 'Sonakshi'
$arr[13,0] = 'Good product but sound not good'
$arr[13,1] = 'Style Name: Google X Series 2024Size: 43'
$arr[13,2] = '4.0 out of 5 stars'
$arr[13,3] = 'sonu jangra'
$arr[14,0] = 'The Xiaomi 55-inch X 4K TV has great picture quality with Dolby Vision and smooth performance. Smart features work well, and the sound is decent. There is some occasional lag, but overall, it''s a good TV for the price.'
$arr[14,1] = 'Style Name: Google X Series 2023Size: 55'
$arr[14,2] = '4.0 out of 5 stars'
$arr[14,3] = $null
$arr[15,0] = 'Very good'
$arr[15,1] = 'Style Name: A Series 2024Size: 43 inches'
$arr[15,2] = '4.0 out of 5 stars'
$arr[15,3] = 'Arslan khan'
$arr[16,0] = 'It’s good for the price but the sound isn’t so great for “Dolby audio” that’s all. Otherwise all good.'
$arr[16,1] = 'Style Name: A Series 2024Size: 43 inches'
$arr[16,2] = '4.0 out of 5 stars'
$arr[16,3] = 'Placeholder'
$arr[17,0] = 'Good product'
$arr[17,1] = 'Style Name: Google X Series 2024Size: 50'
$arr[17,2] = '4.0 out of 5 stars'
$arr[17,3] = 'Amazon buyer'
$arr[18,0] = 'Value for money'
$arr[18,1] = 'Style Name: Google X Series 2024Size: 50'
$arr[18,2] = '4.0 out of 5 stars'
$arr[18,3] = 'Priye Ranjan Kumar'
$arr[19,0] = 'Working well'
$arr[19,1] = 'Style Name: Google X Series 2023Size: 43'
$arr[19,2] = '4.0 out of 5 stars'
$arr[19,3] = 'Pranay'
$arr[20,0] = 'best video quality good picture quality in this price range good performance'
$arr[20,1] = 'Style Name: A Series 2024Size: 43 inches'
$arr[20,2] = '4.0 out of 5 stars'
$arr[20,3] = 'Arvind dabhi'
$arr[21,0] = 'Performance is low, quality is good'
$arr[21,1] = 'Style Name: Google X Series 2024Size: 50'
$arr[21,2] = '4.0 out of 5 stars'
$arr[21,3] = $null
$arr[22,0] = 'Loved It, i kind of believe processing of external media needs improvement, but that can be fixed with software update (hope so),
otherwise, loved the display and audio, speaker are loud, and clear

Quality of remote is not the best, and placement of usb is kinda odd, could have given on bottom

and , the service was good'
$arr[22,1] = 'Style Name: Google X Series 2024Size: 50'
$arr[22,2] = '5.0 out of 5 stars'
$arr[22,3] = 'srinivasa m.'
$arr[23,0] = 'The tv was good working perfectly fine. But the demo of how the tv works wasn''t given properly by the delivery boy,And he delivery boy was arrogant.Overall the product was good but the delivery service wasn''t good.The tv promised the ott apps and other features which were in the tv.'
$arr[23,1] = 'Style Name: A Series 2024Size: 43 inches'
$arr[23,2] = '3.0 out of 5 stars'
$arr[23,3] = 'Baljeet Kaur'
$arr[24,0] = 'Best quality'
$arr[24,1] = 'Style Name: A Series 2024Size: 43 inches'
$arr[24,2] = '4.0 out of 5 stars'
$arr[24,3] = $null
$arr[25,0] = 'I have purchased this TV recently during the amazone sale and delivery and installation is also very fast and good by the team. I have very moderate usage of TV and also TV sound and picture quality is also good as per the price segment.

Recently from last 2 week sometimes I have faced following issues in my TV.

1. When I am turning off my TV by remote then before screen off approx. 50mm green screen is displayed at the bottom side of the TV exact above the TV brand logo area and rest screen will remains normal in color(Black). I don''t know why this happen as this is not all the times happen it is happen sometimes there is no fix frequency or pattern.

2. Sometimes remote voice assistance is not working as not taking any command and all the other function is working normal and only this is not working. This is also not happen routine I have faced this issue 2 times till the purchase. Also I have checked by doing on off the TV but assistance issue is not resolved for couple of minutes and it is resolved within sort time.

Any one have any idea or same issue then pl give your opinion/feedback.'
$arr[25,1] = 'Style Name: Google X Series 2024Size: 55'
$arr[25,2] = '4.0 out of 5 stars'
$arr[25,3] = 'Sridevi'
$arr[26,0] = 'The picture quality and sound quality is superb and the installation service is also good. However, I am not happy with the delivery. The delivery person refused to deliver to the location and behaved very rudely. He was demanding extra money for the delivery.'
$arr[26,1] = 'Style Name: Google X Series 2024Size: 55'
$arr[26,2] = '5.0 out of 5 stars'
$arr[26,3] = 'manoj'
$arr[27,0] = 'Sound balancing is very poor'
$arr[27,1] = 'Style Name: Google X Series 2023Size: 43'
$arr[27,2] = '3.0 out of 5 stars'
$arr[27,3] = 'Bhavik Patel'
$arr[28,0] = 'This tv is outstanding for this panel, picture is good,sound a bit lower at high volume atleast for me but you have a soundbar just connect it and enjoy.
Software also smooth and almost lag free'
$arr[28,1] = 'Style Name: A Series 2024Size: 43 inches'
$arr[28,2] = '4.0 out of 5 stars'
$arr[28,3] = 'Sarwan Kumar'
$arr[29,0] = 'Picture quality, build quality, everything is good except sound quality. Build in speakers are very worst, you have buy a external speaker if want to hear contents clearly'
$arr[29,1] = 'Style Name: A Series 2024Size: 43 inches'
$arr[29,2] = '4.0 out of 5 stars'
$arr[29,3] = '@runsiva'
$arr[30,0] = 'Sound is good
Display is average'
$arr[30,1] = 'Style Name: A Series 2024Size: 43 inches'
$arr[30,2] = '4.0 out of 5 stars'
$arr[30,3] = 'Peddireddy malla Reddy'
$arr[31,0] = 'in this price range this tv is very good. Worth it for buying but there is only one con I noticed in 1 week is Black of this tv not that but no worries its a 27k tv and it works at its best'
$arr[31,1] = 'Style Name: Google X Series 2024Size: 43'
$arr[31,2] = '5.0 out of 5 stars'
$arr[31,3] = 'Placeholder'
$arr[32,0] = 'Good product at low price'
$arr[32,1] = 'Style Name: A Series 2024Size: 43 inches'
$arr[32,2] = '4.0 out of 5 stars'
$arr[32,3] = 'pradeep s'
$arr[33,0] = 'The quality is Ok. A 500 rupees installation fee was charged for the wall mount stand.'
$arr[33,1] = 'Style Name: A Series 2024Size: 43 inches'
$arr[33,2] = '3.0 out of 5 stars'
$arr[33,3] = 'Ravi'
$arr[34,0] = 'I have booked 2gb Ram and 8gb Rom but I have received in this tv only 4Gb Storage Rest everything okay'
$arr[34,1] = 'Style Name: A Series 2024Size: 43 inches'
$arr[34,2] = '4.0 out of 5 stars'
$arr[34,3] = 'narayana av'
$arr[35,0] = 'Click to play video
Xiaomi A Pro 4k 43 inches
Best TV under this budget
I didn’t expected this much Great quality from a budget segment TV.
Colours are very crisp and clear.
Motion smoothening works really well, just switch on it to high and the video quality is going to be very smooth
5 stars go for it You will never regret ??'
$arr[35,1] = 'Style Name: A Series 2024Size: 43 inches'
$arr[35,2] = '5.0 out of 5 stars'
$arr[35,3] = 'Param Saini'
$arr[36,0] = 'TV looks quite decent and good.
Sound was also good but not that Amazing.
Picture clarity is good.

If you need more bass effect definitely you need to plug in External speakers.
UI is good neat clean.

But only thing was the installation part the guy who came to install doesn''t at all look like a professional. Atleast for a company like this they should give contract or appoint with some standard engineer or knowledgeable person.

Rest everything is good. Its being just 1st day I am using it so yeah best of luck to me ✌??..'
$arr[36,1] = 'Style Name: Google X Series 2023Size: 55'
$arr[36,2] = '4.0 out of 5 stars'
$arr[36,3] = $null
$arr[37,0] = 'TV is good in this range. It will become a bit slow but not that much.
Sound is not upto the mark
Picture quality is good.
Installation service was good
I will give it 4 stars.'
$arr[37,1] = 'Style Name: Google X Series 2023Size: 43'
$arr[37,2] = '4.0 out of 5 stars'
$arr[37,3] = 'HV'
$arr[38,0] = 'Ok but it gets stuck in between while watching look n price is very worth'
$arr[38,1] = 'Style Name: Google X Series 2024Size: 55'
$arr[38,2] = '3.0 out of 5 stars'
$arr[38,3] = 'Rahul'
$arr[39,0] = 'Good quality but sound very low'
$arr[39,1] = 'Style Name: Google X Series 2024Size: 43'
$arr[39,2] = '5.0 out of 5 stars'
$arr[39,3] = 'waseem231986'
$arr[40,0] = 'Quality is good and sounds also too good ??'
$arr[40,1] = 'Style Name: A Series 2024Size: 43 inches'
$arr[40,2] = '4.0 out of 5 stars'
$arr[40,3] = 'Pradeep kavali'
$arr[41,0] = 'Picture quality is good but sound quality is really bad.'
$arr[41,1] = 'Style Name: Google X Series 2024Size: 43'
$arr[41,2] = '3.0 out of 5 stars'
$arr[41,3] = $null
$arr[42,0] = 'nyc'
$arr[42,1] = 'Style Name: A Series 2024Size: 43 inches'
$arr[42,2] = '5.0 out of 5 stars'
$arr[42,3] = 'Amazon Customer'
$arr[43,0] = 'This review covers about misleading warranty claim, issues with screen mirroring and latency in TalkBack(screen reading software for visually impaired)
1. The 2 year extended warranty is basically a insurance policy given to us as an OPTION by a insurance company like Zopper and many others. This policy value is capped at 85% of sum of the insured value (which is the tv cost) after depreciation but in the description it says you can claim multiple times as if there is no restrictions at all. This policy becomes null and void once your claim reaches the limit and it doesn''t cover the cost of your claim if it is caused by electric fluctuations and physical dents or damages.
This warranty is absolutely unnecessary since the company itself (Xiaomi) provides 2 year tv panel warranty so, don''t get fooled by this absurd description
2. There is a huge latency when it comes to screen mirroring, it is so much that the video and audio mismatches. It is okay for browsing and just internet surfing but not ideal for entertainment.
3. There are couple of hidden features that is not said in manual or by installation persons, you need to watch a youtube video to know it
4. TalkBack feature is very hard to use since the voice response has latency and generally the tv user interface is not friendly to use for visually impaired persons.'
$arr[43,1] = 'Style Name: Google X Series 2024Size: 55'
$arr[43,2] = '3.0 out of 5 stars'
$arr[43,3] = 'sonam singh'
$arr[44,0] = 'TV IS GOOD BUT REMOTE LAG TOO MUCH'
$arr[44,1] = 'Style Name: A Series 2024Size: 43 inches'
$arr[44,2] = '4.0 out of 5 stars'
$arr[44,3] = 'S.Muruganandam'
$arr[45,0] = 'This tv was delivered to me and I have been using it from two days.

I see the below issues with tv:

1. The picture and sound settings cannot be done over any video. It has to be done seperately in settings and video should be played and tested later. We have to do these trial and error lot of time. (Long press on patchwall button top left button get the setting over video)

2. The internal memory is just 8GB and Ram is 2GB. In that, only 2-3 GB will be remaining after installing basic apps like zee5, sonyliv, hotstar, sunnxt, jio
I got Android update and I am not able to download updates due to space issue.

3. Looks very bright even on standard settings. Tv picture setting is not affecting the setopbox channels and videos. It affects only the videos played via ott apps or YouTube.

4. Going to settop box is not a single button. We need to go to settings-> input source-> HDMI everytime.( Long press on patchwall button and we get option for input source)

5. No direct mute button, we need to press volume down button twice.

There is a lot to explore and will post it later once I discover over coming days.'
$arr[45,1] = 'Style Name: Google X Series 2024Size: 55'
$arr[45,2] = '3.0 out of 5 stars'
$arr[45,3] = 'Mudassir'
$arr[46,0] = 'TV is good, sound & picutre quality superb but lower storage, It have 8GB ROM but 500 MB RAM is only available to use. I try to attached Pendrive , tv is also not support and tb is hanging after attached 128GB pendrive into TV. I am disappoint with low storage.'
$arr[46,1] = 'Style Name: Google X Series 2023Size: 43'
$arr[46,2] = '4.0 out of 5 stars'
$arr[46,3] = 'Arun Ramesha'
$arr[47,0] = 'For product Mentioned Premium "Metal" Bezel Less Design nut there is no metal fo Bezels.'
$arr[47,1] = 'Style Name: Google X Series 2024Size: 43'
$arr[47,2] = '3.0 out of 5 stars'
$arr[47,3] = $null
$arr[48,0] = 'pros:
1. good picture quality
2. good software experience
3. free tv channels support
4. all ports available

cons:
1. sound quality is not able to match the picture quality (Bass missing).
2. takes more time to turn on.
3. Screen mirroring is very laggy
4.sometime tv responds late to remote control.
5. memory is not sufficient for moderate or heavy usage

request/suggestions to Brand:
1. should have brightness control button on remote control.
2. storage can be increased'
$arr[48,1] = 'Style Name: Google X Series 2024Size: 43'
$arr[48,2] = '4.0 out of 5 stars'
$arr[48,3] = 'Manoj'
$arr[49,0] = 'Display is so good but you have to change in settings for your likings, and the sound is the real bummer I had a tv with 20w speakers but they sounded way better than this because it had Dolby atmos, this 30w speakers okish not as expected, the remote is missing a mute feature so my father is having hard time muting the tv because double pressing volume down mutes the tv. And i really suggest you to go for this tv only when you have a home theatre setup or soundbar. And also the viewing angle is not the best you are able see the picture from side but it is greyed out and the brightness is drastically changes when your straight and side so normally tv comes with higher brightness setting which I didn''t like so changed it now it looks good but not from side.'
$arr[49,1] = 'Style Name: A Series 2024Size: 43 inches'
$arr[49,2] = '4.0 out of 5 stars'
$arr[49,3] = 'vaishnavi'
$arr[50,0] = 'Wonderful product'
$arr[50,1] = 'Style Name: Google X Series 2023Size: 43'
$arr[50,2] = '5.0 out of 5 stars'
$arr[50,3] = 'Aritra'
$arr[51,0] = 'Good product but from the price (25k) which i bought this TV, there are better products now in 2024 sadly i was a week early on buying this or else would''ve def not buyed this'
$arr[51,1] = 'Style Name: Google X Series 2023Size: 43'
$arr[51,2] = '4.0 out of 5 stars'
$arr[51,3] = 'VINAYAKAN K S'
$arr[52,0] = 'Nice'
$arr[52,1] = 'Style Name: Google X Series 2024Size: 43'
$arr[52,2] = '5.0 out of 5 stars'
$arr[52,3] = 'Priyanka'
$arr[53,0] = 'Using this product for more than a month. Outstanding product . Value for money. Useless installation service . Third grade installation service .'
$arr[53,1] = 'Style Name: X SeriesSize: 65 inches'
$arr[53,2] = '5.0 out of 5 stars'
$arr[53,3] = 'Ajay Chaudhary'
$arr[54,0] = 'Quality very good
Thank you Xiaomi ??'
$arr[54,1] = 'Style Name: A Series 2024Size: 43 inches'
$arr[54,2] = '5.0 out of 5 stars'
$arr[54,3] = 'GOBINDA PAUL'
$arr[55,0] = 'Good Product at this Price......'
$arr[55,1] = 'Style Name: Google X Series 2023Size: 43'
$arr[55,2] = '5.0 out of 5 stars'
$arr[55,3] = 'Santosh'
$arr[56,0] = 'The worst ever sound quality. Good display but with performance lag! After a week, sound became weaker than it was in the beginning! No more option to replace the TV or replace with another brand!!'
$arr[56,1] = 'Style Name: A Series 2024Size: 43 inches'
$arr[56,2] = '3.0 out of 5 stars'
$arr[56,3] = 'SELVAM'
$arr[57,0] = 'Clarity if fine It has issues of sound Once I shift from OTT platforms to DishTv it becomes silent.Then I have to plug out and plug in couple of times to get the sound.Unable to find service backup'
$arr[57,1] = 'Style Name: A Series 2024Size: 43 inches'
$arr[57,2] = '3.0 out of 5 stars'
$arr[57,3] = 'humayun qaisar'
$arr[58,0] = 'The quality of the TV outstanding but in this price range it would be HDR+ and memory need to improve from 16 GB to 32 for more application support Dolby Atmos not meet as expected overall good ??'
$arr[58,1] = 'Style Name: Google X Series 2024Size: 55'
$arr[58,2] = '4.0 out of 5 stars'
$arr[58,3] = 'MD ASHKAR ALI'
$arr[59,0] = 'Good'
$arr[59,1] = 'Style Name: Google X Series 2023Size: 50'
$arr[59,2] = '4.0 out of 5 stars'
$arr[59,3] = $null
$arr[60,0] = 'This TV lacking of Good Sound Quality with Base, and sound mode'
$arr[60,1] = 'Style Name: A Series 2024Size: 43 inches'
$arr[60,2] = '3.0 out of 5 stars'
$arr[60,3] = 'Sagee'
$arr[61,0] = 'Video quality is very good, but the internal storage of 8gb is not sufficient.
As majority of the space ( about 6gb ) goes of system app and pre installed apps.
I am currently having prime, hotstar, jio. After that I activated Netflix and watched a movie and switched off.
When I tried again to swich on next day it failed. Imagine a 25k rs worth tv not able to take 5 app and failing to get switch on or respond.
After an hour of effort ( as the tv was not responding ) I was able to delete all apps and then installed prime,Netflix and hot star.
that''s all, i cant install one more app.'
$arr[61,1] = 'Style Name: A Series 2024Size: 43 inches'
$arr[61,2] = '3.0 out of 5 stars'
$arr[61,3] = 'Avinash K.'
$arr[62,0] = 'I choosed a wrong platform to order.'
$arr[62,1] = 'Style Name: A Series 2024Size: 43 inches'
$arr[62,2] = '3.0 out of 5 stars'
$arr[62,3] = 'Bishnuyasha Dash'
$arr[63,0] = 'amazing picture and sound quality, but the remote control quality is pathetic'
$arr[63,1] = 'Style Name: Google X Series 2023Size: 43'
$arr[63,2] = '4.0 out of 5 stars'
$arr[63,3] = 'Durga Jyothi Kumar'
$arr[64,0] = 'It has low sound quality'
$arr[64,1] = 'Style Name: A Series 2024Size: 43 inches'
$arr[64,2] = '3.0 out of 5 stars'
$arr[64,3] = 'Placeholder'
$arr[65,0] = 'The bass need to be more in this kind of TV''s I didn''t like the sound of this tv because bass is not that good. The picture quality is fine but I thought if the storage and ram will be higher at this price it would be best. And MI Installation team came 3rd after tv delivery and after calling so many times. The delivery was on date by Amazon but I think it would come much earlier under a weak. One more think this tv doesn''t have patchwall plus I don''t know why this is'
$arr[65,1] = 'Style Name: A Series 2024Size: 43 inches'
$arr[65,2] = '4.0 out of 5 stars'
$arr[65,3] = 'Placeholder'
$arr[66,0] = 'The picture quality of this TV is SERIOUSLY good for the price. Bezels are also minimal, looks great & crisp. Very good colour optimisation. There is no buffering issue as the TV supports 5Ghz Wi-Fi, so streaming 4K HDR content is a breeze if you have a >50-70mbps connection.

UI is smooth for the most part with some lags here & there.

The speakers are awful and not clear at all, would recommend getting a soundbar.

To conclude: you cannot get anything better than this under the ₹20-25K range currently.'
$arr[66,1] = 'Style Name: A Series 2024Size: 43 inches'
$arr[66,2] = '4.0 out of 5 stars'
$arr[66,3] = 'Soumyadipta N.'
$arr[67,0] = 'Stunning picture quality
Dolby Vision
HDR 10 Support
Good sound Dolby Audio Support
Timely Installation
Superb Performance
Wall mount charges ₹499/-'
$arr[67,1] = 'Style Name: A Series 2024Size: 43 inches'
$arr[67,2] = '5.0 out of 5 stars'
$arr[67,3] = $null
$arr[68,0] = 'As i purchase this Mi A Pro tv because of new launch and affordable price.
As i can see some lag in OS as well as playing any 4k content and the internet speed was slow as i connected it through Hotspot...it was taking some time to load even though my internet speed was 300mbps and the delivery guys carelessly kept this Tv flat in their vehicle but tnx there was no damage rather then everything is good worth buying and i hope they fix some bugs in next update
Quick Pros*
Excellent Picture Quality
Excellent Sound
Build quality is good
Affordable Price Worth buying

Cons* ( which i faced )
Little lag in OS & 4k content
Response time is little slow
Stand build quality is not good'
$arr[68,1] = 'Style Name: A Series 2024Size: 43 inches'
$arr[68,2] = '4.0 out of 5 stars'
$arr[68,3] = 'Abhishek C.'
$arr[69,0] = 'Writing this review after a month of purchasing the TV.
Fantastic performance very clear picture and sound.
Fully satisfied with purchasing on Amazon.
Smart TV performance 10 out of 10'
$arr[69,1] = 'Style Name: X SeriesSize: 65 inches'
$arr[69,2] = '5.0 out of 5 stars'
$arr[69,3] = 'Abhishek bagade'
$arr[70,0] = 'Picture quality good best color is so good but voice quality not good'
$arr[70,1] = 'Style Name: A Series 2024Size: 43 inches'
$arr[70,2] = '5.0 out of 5 stars'
$arr[70,3] = 'Shubham Setia'
$arr[71,0] = 'Performance has improved my miles over the previous OLED versions also picture quality and display has improved .
Installation service was also good and the person promptly suggested settings for best viewing
Voice recognition is also improved .
As for features that can be improved , browsing from home is not upto the mark. The TV mounting panel is cheap and it feels that it’s strong enough to support the TV'
$arr[71,1] = 'Style Name: Google X Series 2024Size: 50'
$arr[71,2] = '4.0 out of 5 stars'
$arr[71,3] = 'Debtonoya Das'
$arr[72,0] = 'Audio is good,but after some time like 4 month loading issue'
$arr[72,1] = 'Style Name: Google X Series 2023Size: 43'
$arr[72,2] = '3.0 out of 5 stars'
$arr[72,3] = 'Prakash Kamboya'
$arr[73,0] = 'Very good'
$arr[73,1] = 'Style Name: A Series 2024Size: 43 inches'
$arr[73,2] = '5.0 out of 5 stars'
$arr[73,3] = 'Sanjeev kumar'
$arr[74,0] = 'It was great experience function and quality. Obviously value for money'
$arr[74,1] = 'Style Name: A Series 2024Size: 43 inches'
$arr[74,2] = '4.0 out of 5 stars'
$arr[74,3] = 'Placeholder'
$arr[75,0] = 'GOOD'
$arr[75,1] = 'Style Name: Google X Series 2023Size: 43'
$arr[75,2] = '5.0 out of 5 stars'
$arr[75,3] = 'MANTU SHARMA'
$arr[76,0] = 'I have taken this tv 3 months back, I didn''t even complete my EMI''s for this product but here''s comes the first issue, screen is on loop, by contacting technician I came to know it''s a software issue and it will 800 rupees, thank God I have taken the extra protection plan other wise I have take this to the customer care centre and I should have pay 800. I don''t like to see an issue within 3 months.'
$arr[76,1] = 'Style Name: A Series 2024Size: 43 inches'
$arr[76,2] = '3.0 out of 5 stars'
$arr[76,3] = 'Srikanth goud'
$arr[77,0] = 'TV was delivered by Amazon on time .... But the After Sale Service by Amazon for TV Installation is simply PATHETIC .... I have been a Prime Member of Amazon Since the beginning but never experienced this kind of POOR Service from Amazon ..... No installation Technician came from Amazon for 4 days even after repeated followups and escalations .... 1st they were saying it is Brand installation and then when I asked for the Mi Service ID they told this Amazon Installation ..... Ultimately I had to book a separate installation from Mi and they did it ..... This is WORST experience ever for me from Amazon .... I would recommend that you book the TV installation from the Brand directly instead of Amazon ..... TV itself seems to be fine in this Price Range ....'
$arr[77,1] = 'Style Name: Google X Series 2023Size: 43'
$arr[77,2] = '4.0 out of 5 stars'
$arr[77,3] = 'Souvik Sikdar'
$arr[78,0] = 'Sometime I experience leg but overall good smart TV'
$arr[78,1] = 'Style Name: A Series 2024Size: 43 inches'
$arr[78,2] = '4.0 out of 5 stars'
$arr[78,3] = 'Yatharth Maletiya'
$arr[79,0] = 'Very good Product'
$arr[79,1] = 'Style Name: A Series 2024Size: 43 inches'
$arr[79,2] = '4.0 out of 5 stars'
$arr[79,3] = 'S MAHABOOB BASHA'
$arr[80,0] = 'its good product'
$arr[80,1] = 'Style Name: Google X Series 2024Size: 55'
$arr[80,2] = '4.0 out of 5 stars'
$arr[80,3] = 'Roshan lal Thakur'
$arr[81,0] = 'I am using from last 5 days..and no single lag.. even I am using 2.4 gh WiFi with 40mbps speed. Great working.
Average sound quality, No Bass..so you must connect with sound bar or heavy bluetooth speaker. Overall performance is best. Paisa Vasool TV.'
$arr[81,1] = 'Style Name: A Series 2024Size: 43 inches'
$arr[81,2] = '4.0 out of 5 stars'
$arr[81,3] = 'Prime Customer'
$arr[82,0] = 'Nice product'
$arr[82,1] = 'Style Name: A Series 2024Size: 43 inches'
$arr[82,2] = '5.0 out of 5 stars'
$arr[82,3] = 'Rahul ahirwar'
$arr[83,0] = 'Very Nice ??'
$arr[83,1] = 'Style Name: Google X Series 2024Size: 43'
$arr[83,2] = '5.0 out of 5 stars'
$arr[83,3] = 'RUPDHAR CHHURA'
$arr[84,0] = 'Best colity'
$arr[84,1] = 'Style Name: Google X Series 2023Size: 55'
$arr[84,2] = '5.0 out of 5 stars'
$arr[84,3] = $null
$arr[85,0] = 'Remote control has stopped working after using for 20 days.l tried by changing battery but of no use. It seems remote has some defect.'
$arr[85,1] = 'Style Name: A Series 2024Size: 43 inches'
$arr[85,2] = '4.0 out of 5 stars'
$arr[85,3] = 'Sheetal gupta'
$arr[86,0] = 'Received defective piece on 28th Sep 2024. Totally disappointed ??????.. Installation team from Amazon (TVS) didn''t replace the TV with new one. I had raised the request for technician Visit from brand Company through amazon to check the performance of TV.. Finally Xiaomi replaced the defective TV with new one 16th Oct. Kudoos to Xiaomi company and Mi Installation TV..
About TV.. It''s Superb, Picture quality is awesome. I loved it.. Performance of this is top notch.. ??'
$arr[86,1] = 'Style Name: A Series 2024Size: 43 inches'
$arr[86,2] = '4.0 out of 5 stars'
$arr[86,3] = 'vijay kumar sharma'
$arr[87,0] = 'Colour and contrast ratio is brilliant amazing picture quality and it support HDMI 2.1 which slightly increase more quality in 4k content speaker''s is above average not top notch remote is excellent software is also good based on android 14 which work flawless but some gitter is noticable.'
$arr[87,1] = 'Style Name: A Series 2024Size: 43 inches'
$arr[87,2] = '4.0 out of 5 stars'
$arr[87,3] = 'sikandar'
$arr[88,0] = 'Overall tv is good and budget friendly.
Picture quality is good
Sound is good.'
$arr[88,1] = 'Style Name: Google X Series 2024Size: 43'
$arr[88,2] = '4.0 out of 5 stars'
$arr[88,3] = 'Blue'
$arr[89,0] = 'I didn''t like the sound quality of the TV I expected more from the TV, Treble are on the higher side and bass is very low, but on the other hand it''s good according to the price. For suggestions:- You have to use some soundbars or home theatre for better experience.'
$arr[89,1] = 'Style Name: A Series 2024Size: 43 inches'
$arr[89,2] = '4.0 out of 5 stars'
$arr[89,3] = 'Rishi'
$arr[90,0] = 'I liked the picture quality and other features. However I found the sound quality is average.

I called for technician because we were having issues with WiFi. He came and inspect everything then he told that the issue will persist unless you have excellent signal strength. Means even with good and fair signal strength the apps may hang or work slow.'
$arr[90,1] = 'Style Name: A Series 2024Size: 43 inches'
$arr[90,2] = '4.0 out of 5 stars'
$arr[90,3] = 'Dwarkesh'
$arr[91,0] = 'Remote isn''t working after 2 weeks'
$arr[91,1] = 'Style Name: Google X Series 2023Size: 43'
$arr[91,2] = '3.0 out of 5 stars'
$arr[91,3] = 'joel'
$arr[92,0] = 'Sound quality is not good though it''s dolby.
It''s 4K TV but still not satisfied with the speakers you need to connect external speakers for better experience.
There is not big difference in prices of this (so called affordable tv) and other better tv.
If you can increase your budget by some amt,go for better option I say.'
$arr[92,1] = 'Style Name: A Series 2024Size: 43 inches'
$arr[92,2] = '3.0 out of 5 stars'
$arr[92,3] = 'Shubham Topare'
$arr[93,0] = 'I bought it in August 2024 the packaging was good damage on executive delivered the product at my doorstep and the check the product whether it was damage or not then after 2 hours MI service centre technician called me for installation and they reach to my home within 1 hour and installed the product but I have to pay it extra 500 rupees for the Wall mount because the Wall mount was not including in the box but that is okay because the Wall mount brand is also MI, functionality of the product is good picture quality is good enough features of the product is good but the space of internal storage is not enough(only 8gb,free space is about 4gb)it should be 16GB or more for internal memory, installation was very quick and we are using the TV for last 2 month for me it''s a good product I got it from Amazon at rupees of 20700 and extra 500 I paid for the Wall mount so basically it was a good deal but after that the TV was sold in great Indian festival at rupees of 17700 so I missed the opportunity to buy the product at a low price of that never mind, the product connectivity is good, picture quality was very good just go for it'
$arr[93,1] = 'Style Name: Google X Series 2023Size: 43'
$arr[93,2] = '5.0 out of 5 stars'
$arr[93,3] = 'Akash roy'
$arr[94,0] = 'Just OK OK'
$arr[94,1] = 'Style Name: Google X Series 2023Size: 43'
$arr[94,2] = '4.0 out of 5 stars'
$arr[94,3] = $null
$arr[95,0] = 'Superb Picture Quality and impressive view....nice'
$arr[95,1] = 'Style Name: Google X Series 2023Size: 43'
$arr[95,2] = '4.0 out of 5 stars'
$arr[95,3] = 'MUKESH KUMAR GUPTA'
$arr[96,0] = 'In this range , picture quality sound and other features are very good , value of money product'
$arr[96,1] = 'Style Name: A Series 2024Size: 43 inches'
$arr[96,2] = '4.0 out of 5 stars'
$arr[96,3] = 'MIZANUR RAHMAN'
$arr[97,0] = 'Got my tv in 2 days. Tv display and sound is superb and best in class. Watching youtube videos in 4k 60fps without any lag. However the amazon installation service was pathetic. My time alloted was between 10am - 2pm. But technician didnt came in the time slot, at 4pm evening he called me. So amazon can’t provide installation in the given time slot. Also their technician are liar and cheater. As he asked me 750rs for wall mount stand while it’s clearly mention that stand costs 499rs. My advice would be dont go for amazon installation. Better to ask respective company for installation and demo'
$arr[97,1] = 'Style Name: Google X Series 2024Size: 50'
$arr[97,2] = '4.0 out of 5 stars'
$arr[97,3] = 'BHAGWAN SINGH'
$arr[98,0] = 'Good Led'
$arr[98,1] = 'Style Name: A Series 2024Size: 43 inches'
$arr[98,2] = '5.0 out of 5 stars'
$arr[98,3] = 'Dr. nayak'
$arr[99,0] = 'I purchased this tv due to launch hype and reasonable price, but the tv lags many times due to low available storage,(8gb, out of which 4 gb is available only)if you install any 5 apps out of these (Netflix, hotstar, jiocinema, Sony LIV,zee5, prime video or MX player) , You will definitely feel the lag switching from one to another, also you need to clear cache all the time yo install or update another app, .
Sound is also not upto mark, feels like 300 wala speaker sound. Isse better 32 inch ka LG tv ka sound hai'
$arr[99,1] = 'Style Name: A Series 2024Size: 43 inches'
$arr[99,2] = '3.0 out of 5 stars'
$arr[99,3] = 'Shyam'

$ws.Range("A2:D101").Value = $arr
Write-Host "Done writing data rows"